$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.002045
$ws.Range("H2").Value = 0.12237
$ws.Range("I2").Value = 0.376024
$ws.Range("J2").Value = -1.410655
$ws.Range("K2").Value = -0.824013
$ws.Range("L2").Value = -0.109548
$ws.Range("M2").Value = 1.497178
$ws.Range("N2").Value = -0.696792
$ws.Range("O2").Value = -0.112821
$ws.Range("P2").Value = -0.08856799999999999
$ws.Range("Q2").Value = 1.398434
$ws.Range("R2").Value = -0.153655
$ws.Range("G3").Value = 0.008595
$ws.Range("H3").Value = 0.205328
$ws.Range("I3").Value = 0.361176
$ws.Range("J3").Value = -1.449752
$ws.Range("K3").Value = -0.710593
$ws.Range("L3").Value = -0.091895
$ws.Range("M3").Value = 1.471679
$ws.Range("N3").Value = -0.677548
$ws.Range("O3").Value = -0.096224
$ws.Range("P3").Value = -0.030522
$ws.Range("Q3").Value = 1.182812
$ws.Range("R3").Value = -0.173057
$ws.Range("G4").Value = -0.131351
$ws.Range("H4").Value = -0.007347
$ws.Range("I4").Value = 0.383419
$ws.Range("J4").Value = 1.578679
$ws.Range("K4").Value = -0.02453
$ws.Range("L4").Value = -0.098719
$ws.Range("M4").Value = -0.681806
$ws.Range("N4").Value = 1.155105
$ws.Range("O4").Value = -0.140171
$ws.Range("P4").Value = -0.765522
$ws.Range("Q4").Value = -1.123229
$ws.Range("R4").Value = -0.144528
$ws.Range("G5").Value = -0.136265
$ws.Range("H5").Value = 0.104809
$ws.Range("I5").Value = 0.327248
$ws.Range("J5").Value = 1.56534
$ws.Range("K5").Value = -0.211146
$ws.Range("L5").Value = -0.079442
$ws.Range("M5").Value = -0.927709
$ws.Range("N5").Value = -0.941207
$ws.Range("O5").Value = -0.104856
$ws.Range("P5").Value = -0.501367
$ws.Range("Q5").Value = 1.047544
$ws.Range("R5").Value = -0.14295
$ws.Range("G6").Value = -0.210455
$ws.Range("H6").Value = -0.002195
$ws.Range("I6").Value = 0.335549
$ws.Range("J6").Value = 1.528302
$ws.Range("K6").Value = 0.000402
$ws.Range("L6").Value = -0.06879300000000001
$ws.Range("M6").Value = -0.673835
$ws.Range("N6").Value = 0.900729
$ws.Range("O6").Value = -0.132832
$ws.Range("P6").Value = -0.644012
$ws.Range("Q6").Value = -0.898936
$ws.Range("R6").Value = -0.133923
$ws.Range("G7").Value = -0.00129
$ws.Range("H7").Value = -0.015968
$ws.Range("I7").Value = 0.35939
$ws.Range("J7").Value = -0.756426
$ws.Range("K7").Value = -1.06227
$ws.Range("L7").Value = -0.123168
$ws.Range("M7").Value = -0.535283
$ws.Range("N7").Value = 1.176433
$ws.Range("O7").Value = -0.116348
$ws.Range("P7").Value = 1.292999
$ws.Range("Q7").Value = -0.098195
$ws.Range("R7").Value = -0.119874
$ws.Range("G8").Value = 0.006274
$ws.Range("H8").Value = 0.084373
$ws.Range("I8").Value = 0.316425
$ws.Range("J8").Value = 1.206509
$ws.Range("K8").Value = -0.480704
$ws.Range("L8").Value = -0.098214
$ws.Range("M8").Value = -1.06737
$ws.Range("N8").Value = -0.696082
$ws.Range("O8").Value = -0.091955
$ws.Range("P8").Value = -0.145413
$ws.Range("Q8").Value = 1.092413
$ws.Range("R8").Value = -0.126256
$ws.Range("G9").Value = -0.07199999999999999
$ws.Range("H9").Value = 0.009554999999999999
$ws.Range("I9").Value = 0.302927
$ws.Range("J9").Value = 1.245865
$ws.Range("K9").Value = -0.060324
$ws.Range("L9").Value = -0.085065
$ws.Range("M9").Value = -0.629705
$ws.Range("N9").Value = -0.861822
$ws.Range("O9").Value = -0.107614
$ws.Range("P9").Value = -0.54416
$ws.Range("Q9").Value = 0.912591
$ws.Range("R9").Value = -0.110248
$ws.Range("G10").Value = 0.014031
$ws.Range("H10").Value = 0.008227999999999999
$ws.Range("I10").Value = 0.272333
$ws.Range("J10").Value = 1.04998
$ws.Range("K10").Value = -0.138673
$ws.Range("L10").Value = -0.093754
$ws.Range("M10").Value = -0.630161
$ws.Range("N10").Value = -0.828355
$ws.Range("O10").Value = -0.087287
$ws.Range("P10").Value = -0.433849
$ws.Range("Q10").Value = 0.9588
$ws.Range("R10").Value = -0.091293
